$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.825.52'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '3.691.41'
$ws.Range('E3').Value = '  -0.21%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '647.39'
$ws.Range('E5').Value = '  -4.37%  '
$ws.Range('D6').Value = '161.60'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('D8').Value = '0.501'
$ws.Range('E8').Value = '  +1.07%  '
$ws.Range('E9').Value = '  -1.12%  '
$ws.Range('D10').Value = '7.17'
$ws.Range('E10').Value = '  +0.78%  '
$ws.Range('D11').Value = '0.444'
$ws.Range('E11').Value = '  +0.68%  '
$ws.Range('D12').Value = '0.0000232'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '4.314.81'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = '32.78'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('D15').Value = '3.693.29'
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '69.871.72'
$ws.Range('E16').Value = '  +0.39%  '
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = '16.03'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('D19').Value = '6.53'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').Value = '10.43'
$ws.Range('E20').Value = '  +6.11%  '
$ws.Range('D21').Value = '471.34'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '0.653'
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('D23').Value = '80.02'
$ws.Range('E23').Value = '  -0.70%  '
$ws.Range('D24').Value = '3.840.45'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '0.0000128'
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '10.93'
$ws.Range('E27').Value = '  +0.27%  '
$ws.Range('D28').Value = '9.18'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('D30').Value = '1.71'
$ws.Range('E30').Value = '  -2.22%  '
$ws.Range('D31').Value = '2.02'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('B32').Value = 'Kaspa'
$ws.Range('C32').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D32').Value = '0.167'
$ws.Range('E32').Value = '  +3.20%  '
$ws.Range('B33').Value = 'Binance-PegBSC-USD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').Value = '  -0.21%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '6.54'
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('D35').Value = '26.77'
$ws.Range('E35').Value = '  -0.81%  '
$ws.Range('D36').Value = '3.689.15'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').Value = '8.44'
$ws.Range('E37').Value = '  -0.29%  '
$ws.Range('E38').Value = '  -0.11%  '
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '180.03'
$ws.Range('E39').Value = '  +7.38%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').Value = '5.90'
$ws.Range('E40').Value = '  -5.48%  '
$ws.Range('D41').Value = '2.24'
$ws.Range('E41').Value = '  -0.18%  '
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').Value = '0.0905'
$ws.Range('E43').Value = '  +0.30%  '
$ws.Range('D44').Value = '0.933'
$ws.Range('E44').Value = '  -1.25%  '
$ws.Range('D45').Value = '2.84'
$ws.Range('E45').Value = '  +2.43%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').Value = '47.01'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').Value = '29.32'
$ws.Range('E47').Value = '  +3.45%  '
$ws.Range('D48').Value = '0.000273'
$ws.Range('E48').Value = '  -2.14%  '
$ws.Range('D49').Value = '1.27'
$ws.Range('E49').Value = '  -2.70%  '
$ws.Range('D50').Value = '7.87'
$ws.Range('E50').Value = '  -0.24%  '
$ws.Range('D51').Value = '1.05'
$ws.Range('E51').Value = '  -3.35%  '
